$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds numeric-looking values stored as text. Force text
# formatting before assignment so Excel does not reinterpret the
# string as a number, then restore the default "Normal" style so no
# stray formatting is left behind.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "241.97"
Set-TextValue $ws.Range("D3") "23.40"
Set-TextValue $ws.Range("D4") "5.656"
Set-TextValue $ws.Range("D5") "0.05812"
Set-TextValue $ws.Range("D6") "3.415"
Set-TextValue $ws.Range("D7") "6.467"
Set-TextValue $ws.Range("D8") "1.318"
Set-TextValue $ws.Range("D9") "0.7984"
Set-TextValue $ws.Range("D11") "0.07618"
Set-TextValue $ws.Range("D13") "0.02957"
Set-TextValue $ws.Range("D14") "0.09242"
Set-TextValue $ws.Range("D15") "0.001671"
Set-TextValue $ws.Range("D16") "3.322"
Set-TextValue $ws.Range("D17") "0.04752"
Set-TextValue $ws.Range("D18") "0.0005994"
Set-TextValue $ws.Range("D19") "0.006196"
Set-TextValue $ws.Range("D20") "0.005463"
Set-TextValue $ws.Range("D21") "0.001068"
Set-TextValue $ws.Range("D22") "0.0001501"
Set-TextValue $ws.Range("D23") "3.695"
Set-TextValue $ws.Range("D26") "0.1224"
Set-TextValue $ws.Range("D27") "0.001000"
Set-TextValue $ws.Range("D41") "0.007122"
$ws.Range("E41").Value = "40KickTokenKICK"
$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue $ws.Range("D42") "0.1056"
$ws.Range("E42").Value = "41BKEXTokenBKK"
$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextValue $ws.Range("D43") "0.003173"
$ws.Range("E43").Value = "42CEJICEJIWorstin24h"
Set-TextValue $ws.Range("D44") "0.009528"
$ws.Range("E45").Value = "44ACDXExchangeACXT"
Set-TextValue $ws.Range("D46") "0.00005444"
Set-TextValue $ws.Range("D47") "0.00000000751"
Set-TextValue $ws.Range("D48") "0.7858"
Set-TextValue $ws.Range("D49") "0.1017"
$ws.Range("E49").Value = "48BOLOBOLOBestin24h"
Set-TextValue $ws.Range("D50") "0.00002101"
Set-TextValue $ws.Range("D51") "0.01011"
